# Fruta / hortaliza, semanal
# Inserts one new week's worth of data (3 rows: Especial/Primera/Segunda for
# "Murcott") ahead of the existing "Femacal de La Calera - Mandarina" rows,
# shifting the prior rows (old 529:545) down to (532:548).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 529; everything from 529 downward shifts
# down by three rows (old 529 -> new 532, ... old 545 -> new 548).
$ws.Rows("529:531").Insert()

# Common (constant-across-block) column values for this commodity block.
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$codreg    = 5
$tipo      = "Fruta"
$productoId = 100102
$producto   = "Cítricos"
$categoriaId = 100102004
$categoria   = "Mandarina"
$unidad      = "$/bandeja 10 kilos"
$origen      = "Provincia de Quillota"
$kgUnidad    = 10
$fecha       = 44509

function Set-Row($r, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg) {
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = "Murcott"
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $precioMin
    $ws.Cells.Item($r, 15).Value = $precioMax
    $ws.Cells.Item($r, 16).Value = $precioProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $precioKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}

Set-Row 529 "Especial" 68 6000 6000 6000 600
Set-Row 530 "Primera"  70 5000 5000 5000 500
Set-Row 531 "Segunda"  50 4000 4000 4000 400
